$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.42
$ws.Range("G2").Value = 2.52
$ws.Range("I2").Value = 3.75
$ws.Range("J2").Value = 3.1
$ws.Range("L2").Value = 1.48
$ws.Range("M2").Value = 1.08
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 1.17
$ws.Range("R2").Value = 1.25
$ws.Range("S2").Value = 1.05
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 1.94
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 980
$ws.Range("Y2").Value = 13
$ws.Range("AB2").Value = 980
$ws.Range("AH2").Value = 23
$ws.Range("AJ2").Value = 34

# Row 3
$ws.Range("F3").Value = 1.63
$ws.Range("G3").Value = 1.88
$ws.Range("H3").Value = 1.33
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 4.1
$ws.Range("L3").Value = 1.5
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.11
$ws.Range("P3").Value = 1.25
$ws.Range("Q3").Value = 1.01
$ws.Range("R3").Value = 1.17
$ws.Range("S3").Value = 1.05
$ws.Range("T3").Value = 1.03
$ws.Range("U3").Value = 1.03
$ws.Range("V3").Value = 1.01
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AF3").Value = 1000

# Row 4
$ws.Range("F4").Value = 2.92
$ws.Range("G4").Value = 3.25
$ws.Range("H4").Value = 2.6
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 2.84
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 1.11
$ws.Range("P4").Value = 1.25
$ws.Range("Q4").Value = 1.45
$ws.Range("R4").Value = 1.17
$ws.Range("S4").Value = 4.7
$ws.Range("T4").Value = 1.03
$ws.Range("U4").Value = 1.03
$ws.Range("V4").Value = 1.5
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AM4").Value = 1000

# Row 5
$ws.Range("F5").Value = 1.6
$ws.Range("G5").Value = 6.6
$ws.Range("H5").Value = 1.09
$ws.Range("I5").Value = 980
$ws.Range("J5").Value = 1.2
$ws.Range("K5").Value = 950
$ws.Range("L5").Value = 1.45
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 1.37
$ws.Range("O5").Value = 1.39
$ws.Range("P5").Value = 1.25
$ws.Range("Q5").Value = 1.4
$ws.Range("R5").Value = 1.17
$ws.Range("S5").Value = 2.2
$ws.Range("T5").Value = 1.03
$ws.Range("V5").Value = 1.15
$ws.Range("W5").Value = 1.17
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000

# Row 6
$ws.Range("F6").Value = 1.72
$ws.Range("H6").Value = 5.4
$ws.Range("I6").Value = 6.8
$ws.Range("J6").Value = 3.6
$ws.Range("K6").Value = 3.9
$ws.Range("N6").Value = 1.76
$ws.Range("P6").Value = 1.25
$ws.Range("Q6").Value = 1.38
$ws.Range("V6").Value = 1.18

# Row 7
$ws.Range("G7").Value = 1.12
$ws.Range("J7").Value = 1.03
$ws.Range("K7").Value = 1000
$ws.Range("N7").Value = 3.5
$ws.Range("P7").Value = 3.5
$ws.Range("R7").Value = 1.36
$ws.Range("S7").Value = 1.35
$ws.Range("T7").Value = 1.03
$ws.Range("U7").Value = 1.03
$ws.Range("W7").Value = 1.01
$ws.Range("AF7").Value = 1000
$ws.Range("AN7").Value = 1000

# Row 8
$ws.Range("F8").Value = 1.46
$ws.Range("G8").Value = 1.77
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 2.6
$ws.Range("K8").Value = 500
$ws.Range("L8").Value = 1.45
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 1.11
$ws.Range("O8").Value = 1.39
$ws.Range("P8").Value = 1.25
$ws.Range("Q8").Value = 1.39
$ws.Range("R8").Value = 1.18
$ws.Range("S8").Value = 1.05
$ws.Range("T8").Value = 1.03
$ws.Range("U8").Value = 1.03
$ws.Range("V8").Value = 1.13
$ws.Range("W8").Value = 2.28
$ws.Range("X8").Value = 1000
$ws.Range("Y8").Value = 1000
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 1000
$ws.Range("AC8").Value = 1000
$ws.Range("AD8").Value = 1000
$ws.Range("AE8").Value = 1000
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 1000
$ws.Range("AH8").Value = 1000
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 1000
$ws.Range("AK8").Value = 1000
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 1000

# Row 9
$ws.Range("F9").Value = 1.09
$ws.Range("G9").Value = 1000
$ws.Range("H9").Value = 1.09
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 1.03
$ws.Range("K9").Value = 1000
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 1.02
$ws.Range("P9").Value = 1.25
$ws.Range("Q9").Value = 1.33
$ws.Range("R9").Value = 1.13
$ws.Range("S9").Value = 1.01
$ws.Range("T9").Value = 1.01
$ws.Range("U9").Value = 1.01
$ws.Range("V9").Value = 1.01
$ws.Range("W9").Value = 1.01
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 1000
$ws.Range("AC9").Value = 1000
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000

# Row 10
$ws.Range("F10").Value = 1.2
$ws.Range("G10").Value = 1000
$ws.Range("H10").Value = 1.2
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 1.03
$ws.Range("K10").Value = 1000
$ws.Range("M10").Value = 1.04
$ws.Range("N10").Value = 1.11
$ws.Range("P10").Value = 1.24
$ws.Range("Q10").Value = 1.27
$ws.Range("R10").Value = 1.18
$ws.Range("S10").Value = 1.05
$ws.Range("T10").Value = 1.03
$ws.Range("U10").Value = 1.03
$ws.Range("V10").Value = 1.01
$ws.Range("W10").Value = 1.01
$ws.Range("X10").Value = 1000
$ws.Range("Z10").Value = 1000
$ws.Range("AB10").Value = 1000
$ws.Range("AC10").Value = 1000
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 1000
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000

# Row 11
$ws.Range("G11").Value = 1000
$ws.Range("H11").Value = 1.09
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 2.74
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 1.37
$ws.Range("P11").Value = 1.24
$ws.Range("Q11").Value = 1.42
$ws.Range("R11").Value = 1.18
$ws.Range("S11").Value = 1.05
$ws.Range("T11").Value = 1.03
$ws.Range("U11").Value = 1.03
$ws.Range("V11").Value = 1.26
$ws.Range("W11").Value = 1.25
$ws.Range("Z11").Value = 1000
$ws.Range("AB11").Value = 1000
$ws.Range("AM11").Value = 1000
